# Aggiornamento fino a 27/05
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for the new rows (rows 256 to 269), columns A (date serial), B, C, D
$data = @(
    @(256, 44330, 2, 8,  62.81900274833137),
    @(257, 44331, 2, 8,  62.81900274833137),
    @(258, 44332, 4, 11, 86.37612877895563),
    @(259, 44333, 2, 13, 102.0808794660385),
    @(260, 44334, 1, 13, 102.0808794660385),
    @(261, 44335, 1, 13, 102.0808794660385),
    @(262, 44336, 2, 14, 109.9332548095799),
    @(263, 44337, 3, 15, 117.7856301531213),
    @(264, 44338, 0, 13, 102.0808794660385),
    @(265, 44339, 0, 9,  70.67137809187278),
    @(266, 44340, 2, 9,  70.67137809187278),
    @(267, 44341, 0, 8,  62.81900274833137),
    @(268, 44342, 0, 7,  54.96662740478995),
    @(269, 44343, 0, 5,  39.26187671770711)
)

foreach ($rowData in $data) {
    $r = $rowData[0]

    # Copy formatting (date style) from the cell directly above into the new A cell,
    # then set its value.
    $ws.Range("A" + ($r - 1)).Copy($ws.Range("A" + $r))
    $ws.Cells.Item($r, 1).Value = $rowData[1]

    $ws.Cells.Item($r, 2).Value = $rowData[2]
    $ws.Cells.Item($r, 3).Value = $rowData[3]
    $ws.Cells.Item($r, 4).Value = $rowData[4]
}
